$wb = $excel.ActiveWorkbook

# Both "展览" and "全部类型" sheets carry the same event table and both
# need their "想去人数" (F) / "最低票价" (G) figures refreshed for the
# two affected rows (row 2 and row 5).
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("F2").Value = 325
    $ws.Range("G2").Value = 100

    $ws.Range("F5").Value = 126
}
